$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert two new rows at the correct positions ---
# New listing (ChatGPT / pharmacy-record app) goes in at row 2, pushing the rest down by one.
$ws.Rows.Item(2).Insert()
# New listing (Pokepara auto-like dev) goes in at row 7 (post first insert), pushing the rest down again.
$ws.Rows.Item(7).Insert()

# --- Clear any stale hyperlink definitions before rebuilding them ---
$ws.Hyperlinks.Delete()

# --- Write final cell values for every data row (2-10) ---

$ws.Range("A2").Value = "2025-10-13 18:24:50"
$ws.Range("B2").Value = "【注目】ChatGPTを活用した薬歴アプリ開発の依頼"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5412417"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5412417") | Out-Null
$ws.Range("G2").Value = 398
$ws.Range("H2").Value = "🔥GPT,ChatGPT ◆開発 ◇アプリ"

$ws.Range("A3").Value = "2025-10-13 18:24:50"
$ws.Range("B3").Value = "【GAS開発】配送状況管理の自動化を依頼します"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5412306"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5412306") | Out-Null
$ws.Range("G3").Value = 170
$ws.Range("H3").Value = "◆開発,自動化 ◇管理"

$ws.Range("A4").Value = "2025-10-13 18:24:50"
$ws.Range("B4").Value = "【急募】クリニック向け内視鏡画像システム開発の依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5412233"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5412233") | Out-Null
$ws.Range("G4").Value = 125
$ws.Range("H4").Value = "◆開発,システム開発"

$ws.Range("A5").Value = "2025-10-13 18:24:50"
$ws.Range("B5").Value = "【急募】onedrive上のexcelで自動化システム構築依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5412194"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5412194") | Out-Null
$ws.Range("G5").Value = 95
$ws.Range("H5").Value = "◆自動化"

$ws.Range("A6").Value = "2025-10-13 18:24:50"
$ws.Range("B6").Value = "【急募】スタートアップ向けプロダクト開発のパートナー募集"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5412179"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5412179") | Out-Null
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = "◆開発"

$ws.Range("A7").Value = "2025-10-13 18:24:50"
$ws.Range("B7").Value = "初回 ポケパラの自動いいね等の開発"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5412453"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5412453") | Out-Null
$ws.Range("G7").Value = 63
$ws.Range("H7").Value = "◆開発"

$ws.Range("A8").Value = "2025-10-13 18:24:50"
$ws.Range("B8").Value = "【音声コマンド起動】超小型・低電力レコーダーのプロトタイプ開発"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5412261"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5412261") | Out-Null
$ws.Range("G8").Value = 60
$ws.Range("H8").Value = "◆開発"

$ws.Range("A9").Value = "2025-10-13 18:24:50"
$ws.Range("B9").Value = "微生物の特定と分類を行いたく、画像解析の専門家を探しています!(急いでません!)"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5411887"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5411887") | Out-Null
$ws.Range("G9").Value = 18
$ws.Range("H9").ClearContents()

$ws.Range("A10").Value = "2025-10-13 18:24:50"
$ws.Range("B10").Value = "LINE公式(Lステップ)のリッチメニューの構築"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5412357"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5412357") | Out-Null
$ws.Range("G10").Value = 10
$ws.Range("H10").ClearContents()

# --- Widen column H (skill summary) to fit the new, longer tag text ---
$ws.Columns.Item(8).ColumnWidth = 22.166666666666668
